$d = $word.ActiveDocument

# --- Change 1: collapse the three proofErr-wrapped runs ("0 ", " ", "based")
#     into a clean "0 based" run (removing the now-unneeded gramStart/gramEnd
#     markers along with the mid-sentence split).
$d.Content.Find.Execute("0  based", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "0 based", 2)

# --- Change 2: append a new bullet after the "This keyword..." paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*keyword will work depend on how the function have been called.*") {
        $target = $p
    }
}

$newPara = $target.Range.Paragraphs.Add($target.Range)
$newPara.Range.Text = "Whenever you are returning another function make sure you use bind so everything will work just fine."
